$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 23.24000000000019
$ws.Range("G2").Value = 0.0002068047250504135
$ws.Range("H2").Value = 0.001000294255872531
$ws.Range("K2").Value = 4.65237144908373
$ws.Range("L2").Value = "[1.6967320139249722, 7.608010884242487]"
$ws.Range("M2").Value = 0.002105026404423338
$ws.Range("N2").Value = 0.002105026404423338
$ws.Range("O2").Value = -0.9182633181663862
$ws.Range("P2").Value = "[-1.534631846798618, -0.3018947895341544]"
$ws.Range("Q2").Value = 0.003590793943548665
$ws.Range("R2").Value = 0.003590793943548665
$ws.Range("S2").Value = 13.73632025622117
$ws.Range("T2").Value = "[12.169740125773828, 15.302900386668519]"
$ws.Range("W2").Value = 3.396436436436467
$ws.Range("X2").Value = 1.116636636636648
$ws.Range("Y2").Value = 5.676236236236286

$ws.Range("E3").Value = 23.79000000000028
$ws.Range("G3").Value = 0.0000103526522280406
$ws.Range("H3").Value = 0.0002408790199947342
$ws.Range("K3").Value = 4.861689453340309
$ws.Range("L3").Value = "[2.534248141882461, 7.189130764798156]"
$ws.Range("M3").Value = 0.00004766166157366136
$ws.Range("N3").Value = 0.00009532332314732272
$ws.Range("O3").Value = 1.478026573760963
$ws.Range("P3").Value = "[0.9119738433844233, 2.0440793041375027]"
$ws.Range("Q3").Value = 0.0000004207644115350462
$ws.Range("R3").Value = 0.0000008415288230700924
$ws.Range("S3").Value = 13.06188950652026
$ws.Range("T3").Value = "[11.672658469611907, 14.451120543428612]"
$ws.Range("W3").Value = 18.19375375375397
$ws.Range("X3").Value = 16.0505105105107
$ws.Range("Y3").Value = 20.33699699699724

